# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (want-to-go count) figures in column F across the
# four worksheets of the 漫展信息 workbook.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 685
$ws.Cells.Item(6, 6).Value = 113
$ws.Cells.Item(7, 6).Value = 1158
$ws.Cells.Item(10, 6).Value = 2060
$ws.Cells.Item(11, 6).Value = 57
$ws.Cells.Item(12, 6).Value = 40
$ws.Cells.Item(16, 6).Value = 1478
$ws.Cells.Item(17, 6).Value = 1478
$ws.Cells.Item(19, 6).Value = 557
$ws.Cells.Item(20, 6).Value = 395
$ws.Cells.Item(21, 6).Value = 395
$ws.Cells.Item(22, 6).Value = 738
$ws.Cells.Item(23, 6).Value = 451
$ws.Cells.Item(24, 6).Value = 2865
$ws.Cells.Item(25, 6).Value = 400
$ws.Cells.Item(27, 6).Value = 3211
$ws.Cells.Item(28, 6).Value = 663
$ws.Cells.Item(29, 6).Value = 532
$ws.Cells.Item(30, 6).Value = 240
$ws.Cells.Item(31, 6).Value = 983
$ws.Cells.Item(32, 6).Value = 735
$ws.Cells.Item(34, 6).Value = 711
$ws.Cells.Item(35, 6).Value = 693

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(13, 6).Value = 80
$ws.Cells.Item(20, 6).Value = 96
$ws.Cells.Item(21, 6).Value = 190
$ws.Cells.Item(22, 6).Value = 134
$ws.Cells.Item(23, 6).Value = 448

# 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(6, 6).Value = 395

# 全部类型 (All Types - combined listing)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(8, 6).Value = 685
$ws.Cells.Item(10, 6).Value = 113
$ws.Cells.Item(12, 6).Value = 1158
$ws.Cells.Item(15, 6).Value = 395
$ws.Cells.Item(16, 6).Value = 2060
$ws.Cells.Item(17, 6).Value = 57
$ws.Cells.Item(18, 6).Value = 40
$ws.Cells.Item(25, 6).Value = 80
$ws.Cells.Item(26, 6).Value = 1478
$ws.Cells.Item(27, 6).Value = 1478
$ws.Cells.Item(30, 6).Value = 557
$ws.Cells.Item(31, 6).Value = 395
$ws.Cells.Item(32, 6).Value = 395
$ws.Cells.Item(34, 6).Value = 738
$ws.Cells.Item(35, 6).Value = 451
$ws.Cells.Item(37, 6).Value = 2865
$ws.Cells.Item(39, 6).Value = 3211
$ws.Cells.Item(40, 6).Value = 663
$ws.Cells.Item(41, 6).Value = 532
$ws.Cells.Item(42, 6).Value = 240
$ws.Cells.Item(43, 6).Value = 983
$ws.Cells.Item(45, 6).Value = 96
$ws.Cells.Item(46, 6).Value = 190
$ws.Cells.Item(47, 6).Value = 134
$ws.Cells.Item(48, 6).Value = 448
$ws.Cells.Item(49, 6).Value = 735
$ws.Cells.Item(50, 6).Value = 711
$ws.Cells.Item(51, 6).Value = 693
